$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full "boxed header" format (bold font, thin border all around,
# centered horizontally, top vertically) on A2 first ...
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("A2").Borders.Weight = 2
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160

# ... then clone that exact format onto B1 via copy/paste-special so both
# cells end up sharing a single cell style record.
$ws.Range("A2").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = 0

# Plain text label in B2 (default formatting, shared string table entry).
$ws.Range("B2").Value = "disconnected_elements"
